# Scheduled market-price refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H, I, J, K, L, M, N) on the affected leve rows across
# all eight job sheets, per the latest Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 534
$ws.Range("J55").Value = 757.8
$ws.Range("L55").Value = 757.8
$ws.Range("N55").Value = -1185.8
$ws.Range("H69").Value = 11877.777
$ws.Range("J69").Value = 7950
$ws.Range("L69").Value = 23850
$ws.Range("N69").Value = -25598
$ws.Range("H72").Value = 11877.777
$ws.Range("J72").Value = 7950
$ws.Range("L72").Value = 71550
$ws.Range("N72").Value = -80286
$ws.Range("H107").Value = 406.1111
$ws.Range("I107").Value = 423.07693
$ws.Range("K107").Value = 423.07693
$ws.Range("M107").Value = 1496.92307
$ws.Range("H137").Value = 3637.0588
$ws.Range("I137").Value = 1589.6666
$ws.Range("K137").Value = 4768.9998
$ws.Range("M137").Value = -2218.9998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1527.5
$ws.Range("I2").Value = 1505
$ws.Range("K2").Value = 1505
$ws.Range("M2").Value = -1392
$ws.Range("H32").Value = 16110.549
$ws.Range("I32").Value = 7185.839
$ws.Range("K32").Value = 7185.839
$ws.Range("M32").Value = -6898.839
$ws.Range("H45").Value = 2141.4285
$ws.Range("I45").Value = 2165
$ws.Range("K45").Value = 2165
$ws.Range("M45").Value = -1788
$ws.Range("H74").Value = 2247.5625
$ws.Range("I74").Value = 1075.4445
$ws.Range("K74").Value = 1075.4445
$ws.Range("M74").Value = -201.4445000000001
$ws.Range("H77").Value = 2247.5625
$ws.Range("I77").Value = 1075.4445
$ws.Range("K77").Value = 5377.2225
$ws.Range("M77").Value = -1009.2225
$ws.Range("H81").Value = 45000
$ws.Range("I81").Value = 40000
$ws.Range("K81").Value = 40000
$ws.Range("M81").Value = -39002
$ws.Range("H84").Value = 45000
$ws.Range("I84").Value = 40000
$ws.Range("K84").Value = 120000
$ws.Range("M84").Value = -115008
$ws.Range("H116").Value = 1527.5
$ws.Range("I116").Value = 1505
$ws.Range("K116").Value = 1505
$ws.Range("M116").Value = 789
$ws.Range("H122").Value = 4567.25
$ws.Range("I122").Value = 4303.3335
$ws.Range("K122").Value = 12910.0005
$ws.Range("M122").Value = -10460.0005
$ws.Range("H132").Value = 1659.7894
$ws.Range("I132").Value = 1602.2075
$ws.Range("K132").Value = 4806.622499999999
$ws.Range("M132").Value = -2276.622499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1527.5
$ws.Range("I3").Value = 1505
$ws.Range("K3").Value = 1505
$ws.Range("M3").Value = -1391
$ws.Range("H134").Value = 888.24243
$ws.Range("I134").Value = 728.53125
$ws.Range("K134").Value = 2185.59375
$ws.Range("M134").Value = 349.40625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5480.5713
$ws.Range("J31").Value = 7482.375
$ws.Range("L31").Value = 7482.375
$ws.Range("N31").Value = -8072.375
$ws.Range("H34").Value = 5480.5713
$ws.Range("J34").Value = 7482.375
$ws.Range("L34").Value = 7482.375
$ws.Range("N34").Value = -7886.375
$ws.Range("H107").Value = 448.22223
$ws.Range("I107").Value = 317.46155
$ws.Range("K107").Value = 317.46155
$ws.Range("M107").Value = 1602.53845
$ws.Range("H132").Value = 2652.9355
$ws.Range("I132").Value = 2408.1304
$ws.Range("K132").Value = 7224.3912
$ws.Range("M132").Value = -4694.3912

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 400
$ws.Range("I21").Value = 400
$ws.Range("K21").Value = 1200
$ws.Range("M21").Value = -1027
$ws.Range("H75").Value = 701.125
$ws.Range("J75").Value = 772.4
$ws.Range("L75").Value = 2317.2
$ws.Range("N75").Value = -4313.2
$ws.Range("H78").Value = 701.125
$ws.Range("J78").Value = 772.4
$ws.Range("L78").Value = 6951.599999999999
$ws.Range("N78").Value = -16935.6
$ws.Range("H107").Value = 1261.5714
$ws.Range("I107").Value = 1271.5
$ws.Range("K107").Value = 3814.5
$ws.Range("M107").Value = -1894.5
$ws.Range("H113").Value = 1631.75
$ws.Range("I113").Value = 1050
$ws.Range("J113").Value = 1748.1
$ws.Range("K113").Value = 3150
$ws.Range("L113").Value = 5244.299999999999
$ws.Range("M113").Value = -980
$ws.Range("N113").Value = -9584.299999999999
$ws.Range("H140").Value = 2281.4119
$ws.Range("I140").Value = 2281.4119
$ws.Range("K140").Value = 6844.2357
$ws.Range("M140").Value = -1664.2357

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 581451.5
$ws.Range("J122").Value = 2004797.2
$ws.Range("L122").Value = 6014391.6
$ws.Range("N122").Value = -6019291.6
$ws.Range("H126").Value = 3774.9
$ws.Range("I126").Value = 2749.6667
$ws.Range("J126").Value = 4214.2856
$ws.Range("K126").Value = 8249.000100000001
$ws.Range("L126").Value = 12642.8568
$ws.Range("M126").Value = -5779.000100000001
$ws.Range("N126").Value = -17582.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3166.3333
$ws.Range("J40").Value = 3330.6667
$ws.Range("L40").Value = 3330.6667
$ws.Range("N40").Value = -3602.6667
$ws.Range("H61").Value = 5458.3335
$ws.Range("I61").Value = 5749.8
$ws.Range("K61").Value = 5749.8
$ws.Range("M61").Value = -5547.8
$ws.Range("H74").Value = 42999
$ws.Range("I74").Value = 42999
$ws.Range("K74").Value = 42999
$ws.Range("M74").Value = -42001
$ws.Range("H77").Value = 42999
$ws.Range("I77").Value = 42999
$ws.Range("K77").Value = 128997
$ws.Range("M77").Value = -124005
$ws.Range("H100").Value = 1298.75
$ws.Range("I100").Value = 1298.75
$ws.Range("K100").Value = 1298.75
$ws.Range("M100").Value = -757.75
$ws.Range("H113").Value = 5458.3335
$ws.Range("I113").Value = 5749.8
$ws.Range("K113").Value = 5749.8
$ws.Range("M113").Value = -3579.8
$ws.Range("H132").Value = 3018.3225
$ws.Range("I132").Value = 2052.9565
$ws.Range("J132").Value = 5793.75
$ws.Range("K132").Value = 6158.869499999999
$ws.Range("L132").Value = 17381.25
$ws.Range("M132").Value = -3628.869499999999
$ws.Range("N132").Value = -22441.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7662.467
$ws.Range("I62").Value = 6784.25
$ws.Range("J62").Value = 7981.8184
$ws.Range("K62").Value = 6784.25
$ws.Range("L62").Value = 7981.8184
$ws.Range("M62").Value = -6160.25
$ws.Range("N62").Value = -9229.8184
$ws.Range("H65").Value = 7662.467
$ws.Range("I65").Value = 6784.25
$ws.Range("J65").Value = 7981.8184
$ws.Range("K65").Value = 33921.25
$ws.Range("L65").Value = 39909.092
$ws.Range("M65").Value = -30801.25
$ws.Range("N65").Value = -46149.092
$ws.Range("H107").Value = 1475.1666
$ws.Range("I107").Value = 662.5
$ws.Range("K107").Value = 1987.5
$ws.Range("M107").Value = -67.5
$ws.Range("H135").Value = 79249.5
$ws.Range("J135").Value = 79249.5
$ws.Range("L135").Value = 79249.5
$ws.Range("N135").Value = -89389.5
